$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay TEXT (not be coerced to a
# number). Force text format first, then write the value, then restore the
# original cell style (PasteSpecial Formats keeps the s= index intact,
# unlike re-assigning .Style which was losing the format).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B8").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C3").Value = "Mohaupt"
$ws.Range("D5").Value = "KONTOSTAND AM 12.02.2025"

# --- Row 6 ---
$ws.Range("B6").Value = "13.02."
$ws.Range("C6").Value = "14.02."
$ws.Range("D6").Value = "PAYPAL JVJLYJ"
$ws.Range("E6").Value = "7,83-"

# --- Row 7 ---
$ws.Range("B7").Value = "16.02."
$ws.Range("C7").Value = "17.02."
$ws.Range("D7").Value = "BURGER KING Coburg"
$ws.Range("E7").Value = "32,05-"

# --- Row 8 ---
$ws.Range("B8").Value = "17.02."
$ws.Range("C8").Value = "18.02."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,72-"

# --- Row 9 (previously an empty placeholder row) ---
$ws.Range("B9").Value = "18.02."
$ws.Range("C9").Value = "19.02."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 84792294"
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)  # xlPasteFormats (matches target s="17")
$ws.Range("E9").Value = "41,95-"

# --- Row 10 (previously an empty placeholder row) ---
$ws.Range("B10").Value = "22.02."
$ws.Range("C10").Value = "23.02."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-97367660"
$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats (matches target s="17")
$ws.Range("E10").Value = "52,99-"

# --- Row 11 (previously an empty placeholder row) ---
$ws.Range("B11").Value = "24.02."
$ws.Range("C11").Value = "25.02."
$ws.Range("D11").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 77633790"
$ws.Range("E8").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats (matches target s="17")
$ws.Range("E11").Value = "86,64-"

# --- Footer fields ---
$ws.Range("D12").Value = "KONTOSTAND AM 01.03.2025"
$ws.Range("E12").Value = "246,18-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.03.2025"
